$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 18, shifting existing rows 18:48 down to 19:49.
$ws.Rows("18:18").Insert()

# Populate the newly inserted row 18 with the new record.
$ws.Range("A18").Value = 5
$ws.Range("B18").Value = "Macroferia Regional de Talca"
$ws.Range("C18").Value = "Maule"
$ws.Range("D18").Value = 44498
$ws.Range("E18").Value = 7
$ws.Range("F18").Value = 100112026
$ws.Range("G18").Value = "Haba"
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 500
$ws.Range("K18").Value = 6000
$ws.Range("L18").Value = 6000
$ws.Range("M18").Value = 6000
$ws.Range("N18").Value = '$/saco 25 kilos'
$ws.Range("O18").Value = "Región del Maule"
$ws.Range("P18").Value = 240
$ws.Range("Q18").Value = 25
$ws.Range("R18").Value = "Hortaliza"
